$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1) Row 6 (Arm_01): update damage/bulletSpeed/fireRate/coolDownTime/projectilesPerShot/shotErrorRange ---
$ws.Range("O6").Value = 5
$ws.Range("P6").Value = 25
$ws.Range("Q6").Value = 0.1
$ws.Range("R6").Value = 3
$ws.Range("S6").Value = 1
$ws.Range("T6").Value = 2

# --- 2) Push the divider/blank formatting rows down by two, without disturbing rows 12-19 ---
# Move old row 20 (bottom border divider row) down to row 22.
$ws.Rows.Item(20).Copy()
$ws.Rows.Item(22).PasteSpecial()
# Rows 20 and 21 become new blank formatted rows, matching the style of row 19.
$ws.Rows.Item(19).Copy()
$ws.Rows.Item(20).PasteSpecial()
$ws.Rows.Item(19).Copy()
$ws.Rows.Item(21).PasteSpecial()

# --- 3) Copy existing rows 8 & 9 (Shoulder_01 / Shoulder_02 data) down into the now-free rows 10 & 11 ---
$ws.Rows.Item(8).Copy()
$ws.Rows.Item(10).PasteSpecial()
$ws.Rows.Item(9).Copy()
$ws.Rows.Item(11).PasteSpecial()

# --- 4) Overwrite rows 8 & 9 with the new Arm_03 / Arm_04 weapon-arm data ---
$ws.Range("A8").Value = 10003003
$ws.Range("B8").Value = "Arm_03"
$ws.Range("C8").Value = "Prefabs/Parts/Weapon_Arm/Weapon_Arm_03"
$ws.Range("D8").Value = "TEST WEAPON_ARM_03"
$ws.Range("E8").Value = "교체 여부 확인용 파츠"
$ws.Range("F8").Value = "Weapon_Arm"
$ws.Range("G8").Value = 275
$ws.Range("H8").Value = 10
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = $false
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = "Prefabs/Projectiles/PlayerBullet/Bullet_Gatling_01"
$ws.Range("O8").Value = 2
$ws.Range("P8").Value = 3
$ws.Range("Q8").Value = 4
$ws.Range("R8").Value = 5
$ws.Range("S8").Value = 6
$ws.Range("T8").Value = 7

$ws.Range("A9").Value = 10003004
$ws.Range("B9").Value = "Arm_04"
$ws.Range("C9").Value = "Prefabs/Parts/Weapon_Arm/Weapon_Arm_04"
$ws.Range("D9").Value = "TEST WEAPON_ARM_04"
$ws.Range("E9").Value = "교체 여부 확인용 파츠"
$ws.Range("F9").Value = "Weapon_Arm"
$ws.Range("G9").Value = 275
$ws.Range("H9").Value = 10
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = $false
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = "Prefabs/Projectiles/PlayerBullet/Bullet_Gatling_01"
$ws.Range("O9").Value = 2
$ws.Range("P9").Value = 3
$ws.Range("Q9").Value = 4
$ws.Range("R9").Value = 5
$ws.Range("S9").Value = 6
$ws.Range("T9").Value = 7

# --- 5) The copied-down Shoulder_01 row (now row 10) gets a corrected missile bullet prefab
#         plus new combat stats ---
$ws.Range("N10").Value = "Prefabs/Projectiles/PlayerBullet/Bullet_Missile_01"
$ws.Range("O10").Value = 5
$ws.Range("P10").Value = 15
$ws.Range("Q10").Value = 0.05
$ws.Range("R10").Value = 5
$ws.Range("S10").Value = 5
$ws.Range("T10").Value = 1

# --- 6) Restore the sheet selection to match the final editing position ---
$ws.Range("Q7").Select()
